$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 127450
$ws.Range("I106").Value = 168433.33
$ws.Range("K106").Value = 168433.33
$ws.Range("M106").Value = -167802.33
$ws.Range("H112").Value = 1787.2354
$ws.Range("J112").Value = 1845.8125
$ws.Range("L112").Value = 5537.4375
$ws.Range("N112").Value = -7753.4375
$ws.Range("H113").Value = 2404.4443
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 2728
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 2728
$ws.Range("M113").Value = 1254
$ws.Range("N113").Value = -9236
$ws.Range("H132").Value = 4312.643
$ws.Range("I132").Value = 4112.75
$ws.Range("J132").Value = 4812.375
$ws.Range("K132").Value = 12338.25
$ws.Range("L132").Value = 14437.125
$ws.Range("M132").Value = -9808.25
$ws.Range("N132").Value = -19497.125
$ws.Range("H137").Value = 13335074
$ws.Range("I137").Value = 1125.75
$ws.Range("J137").Value = 28573872
$ws.Range("K137").Value = 3377.25
$ws.Range("L137").Value = 85721616
$ws.Range("M137").Value = -827.25
$ws.Range("N137").Value = -85726716
$ws.Range("H138").Value = 2397.4524
$ws.Range("I138").Value = 1787.1111
$ws.Range("J138").Value = 2855.2083
$ws.Range("K138").Value = 5361.3333
$ws.Range("L138").Value = 8565.624899999999
$ws.Range("M138").Value = -221.3333000000002
$ws.Range("N138").Value = -18845.6249
$ws.Range("H141").Value = 4817.1816
$ws.Range("I141").Value = 3997.5
$ws.Range("J141").Value = 4999.3335
$ws.Range("K141").Value = 11992.5
$ws.Range("L141").Value = 14998.0005
$ws.Range("M141").Value = -6812.5
$ws.Range("N141").Value = -25358.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 17243692
$ws.Range("I61").Value = 18520854
$ws.Range("K61").Value = 18520854
$ws.Range("M61").Value = -18520642
$ws.Range("H74").Value = 25004988
$ws.Range("I74").Value = 35717630
$ws.Range("J74").Value = 8821.333000000001
$ws.Range("K74").Value = 35717630
$ws.Range("L74").Value = 8821.333000000001
$ws.Range("M74").Value = -35716756
$ws.Range("N74").Value = -10569.333
$ws.Range("H77").Value = 25004988
$ws.Range("I77").Value = 35717630
$ws.Range("J77").Value = 8821.333000000001
$ws.Range("K77").Value = 178588150
$ws.Range("L77").Value = 44106.665
$ws.Range("M77").Value = -178583782
$ws.Range("N77").Value = -52842.665
$ws.Range("H124").Value = 39250
$ws.Range("J124").Value = 39250
$ws.Range("L124").Value = 39250
$ws.Range("N124").Value = -49070
$ws.Range("H136").Value = 17243692
$ws.Range("I136").Value = 18520854
$ws.Range("K136").Value = 55562562
$ws.Range("M136").Value = -55560012
$ws.Range("H141").Value = 45675
$ws.Range("J141").Value = 45675
$ws.Range("L141").Value = 45675
$ws.Range("N141").Value = -56035

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 50000
$ws.Range("J81").Value = 50000
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -52122
$ws.Range("H84").Value = 50000
$ws.Range("J84").Value = 50000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -160608
$ws.Range("H134").Value = 3496.55
$ws.Range("I134").Value = 2117.9312
$ws.Range("J134").Value = 7131.091
$ws.Range("K134").Value = 6353.7936
$ws.Range("L134").Value = 21393.273
$ws.Range("M134").Value = -3818.7936
$ws.Range("N134").Value = -26463.273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5851339.5
$ws.Range("I31").Value = 3809.8538
$ws.Range("J31").Value = 20835634
$ws.Range("K31").Value = 3809.8538
$ws.Range("L31").Value = 20835634
$ws.Range("M31").Value = -3514.8538
$ws.Range("N31").Value = -20836224
$ws.Range("H34").Value = 5851339.5
$ws.Range("I34").Value = 3809.8538
$ws.Range("J34").Value = 20835634
$ws.Range("K34").Value = 3809.8538
$ws.Range("L34").Value = 20835634
$ws.Range("M34").Value = -3607.8538
$ws.Range("N34").Value = -20836038
$ws.Range("H105").Value = 1421.1177
$ws.Range("I105").Value = 1506
$ws.Range("J105").Value = 1299.8572
$ws.Range("K105").Value = 1506
$ws.Range("L105").Value = 1299.8572
$ws.Range("M105").Value = 241
$ws.Range("N105").Value = -4793.8572
$ws.Range("H107").Value = 1041.4
$ws.Range("I107").Value = 694.2857
$ws.Range("J107").Value = 1851.3334
$ws.Range("K107").Value = 694.2857
$ws.Range("L107").Value = 1851.3334
$ws.Range("M107").Value = 1225.7143
$ws.Range("N107").Value = -5691.3334
$ws.Range("H132").Value = 13160031
$ws.Range("I132").Value = 17859076
$ws.Range("J132").Value = 2703.7
$ws.Range("K132").Value = 53577228
$ws.Range("L132").Value = 8111.099999999999
$ws.Range("M132").Value = -53574698
$ws.Range("N132").Value = -13171.1
$ws.Range("H134").Value = 701928
$ws.Range("I134").Value = 1610.2593
$ws.Range("J134").Value = 3403153.5
$ws.Range("K134").Value = 4830.7779
$ws.Range("L134").Value = 10209460.5
$ws.Range("M134").Value = -2295.7779
$ws.Range("N134").Value = -10214530.5
$ws.Range("H140").Value = 41012.5
$ws.Range("J140").Value = 41012.5
$ws.Range("L140").Value = 41012.5
$ws.Range("N140").Value = -51372.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1117.3158
$ws.Range("I5").Value = 357.18182
$ws.Range("J5").Value = 2162.5
$ws.Range("K5").Value = 1071.54546
$ws.Range("L5").Value = 6487.5
$ws.Range("M5").Value = -959.54546
$ws.Range("N5").Value = -6711.5
$ws.Range("H69").Value = 1342.5264
$ws.Range("I69").Value = 750.5714
$ws.Range("J69").Value = 3000
$ws.Range("K69").Value = 2251.7142
$ws.Range("L69").Value = 9000
$ws.Range("M69").Value = -1440.7142
$ws.Range("N69").Value = -10622
$ws.Range("H72").Value = 1342.5264
$ws.Range("I72").Value = 750.5714
$ws.Range("J72").Value = 3000
$ws.Range("K72").Value = 6755.1426
$ws.Range("L72").Value = 27000
$ws.Range("M72").Value = -2699.1426
$ws.Range("N72").Value = -35112
$ws.Range("H131").Value = 1215
$ws.Range("I131").Value = 473.84616
$ws.Range("J131").Value = 1420
$ws.Range("K131").Value = 1421.53848
$ws.Range("L131").Value = 4260
$ws.Range("M131").Value = 3618.46152
$ws.Range("N131").Value = -14340
$ws.Range("H132").Value = 1215.4333
$ws.Range("I132").Value = 664.1429000000001
$ws.Range("J132").Value = 1697.8125
$ws.Range("K132").Value = 5977.2861
$ws.Range("L132").Value = 15280.3125
$ws.Range("M132").Value = -3447.2861
$ws.Range("N132").Value = -20340.3125
$ws.Range("H133").Value = 111115400
$ws.Range("I133").Value = 222224200
$ws.Range("J133").Value = 6610
$ws.Range("K133").Value = 666672600
$ws.Range("L133").Value = 19830
$ws.Range("M133").Value = -666667540
$ws.Range("N133").Value = -29950
$ws.Range("H135").Value = 1117.3158
$ws.Range("I135").Value = 357.18182
$ws.Range("J135").Value = 2162.5
$ws.Range("K135").Value = 3214.63638
$ws.Range("L135").Value = 19462.5
$ws.Range("M135").Value = -679.6363799999999
$ws.Range("N135").Value = -24532.5
$ws.Range("H141").Value = 10348.272
$ws.Range("I141").Value = 7978.875
$ws.Range("K141").Value = 23936.625
$ws.Range("M141").Value = -18756.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 38
$ws.Range("I2").Value = 40
$ws.Range("J2").Value = 35.5
$ws.Range("K2").Value = 40
$ws.Range("L2").Value = 35.5
$ws.Range("M2").Value = 73
$ws.Range("N2").Value = -261.5
$ws.Range("H127").Value = 27666.666
$ws.Range("J127").Value = 27666.666
$ws.Range("L127").Value = 27666.666
$ws.Range("N127").Value = -37586.666
$ws.Range("H132").Value = 4209.2856
$ws.Range("I132").Value = 4039.8667
$ws.Range("J132").Value = 4632.8335
$ws.Range("K132").Value = 12119.6001
$ws.Range("L132").Value = 13898.5005
$ws.Range("M132").Value = -9589.6001
$ws.Range("N132").Value = -18958.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7204.7856
$ws.Range("I40").Value = 8107.4443
$ws.Range("J40").Value = 5580
$ws.Range("K40").Value = 8107.4443
$ws.Range("L40").Value = 5580
$ws.Range("M40").Value = -7971.4443
$ws.Range("N40").Value = -5852

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1143.1428
$ws.Range("I136").Value = 1027.92
$ws.Range("J136").Value = 2103.3333
$ws.Range("K136").Value = 3083.76
$ws.Range("L136").Value = 6309.999899999999
$ws.Range("M136").Value = -533.7600000000002
$ws.Range("N136").Value = -11409.9999
